$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" '37.162.37'
$ws.Range("E2").Value = '  +1.79%  '
Set-TextValue "D3" '2.019.76'
$ws.Range("E3").Value = '  +3.32%  '
$ws.Range("E4").Value = '  +0.05%  '
Set-TextValue "D5" '246.61'
$ws.Range("E5").Value = '  +1.41%  '
Set-TextValue "D6" '0.626'
$ws.Range("E6").Value = '  -0.13%  '
Set-TextValue "D7" '60.23'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("E9").Value = '  +3.02%  '
Set-TextValue "D10" '0.0807'
$ws.Range("E10").Value = '  +2.21%  '
$ws.Range("E11").Value = '  +1.27%  '
Set-TextValue "D12" '14.98'
$ws.Range("E12").Value = '  +5.86%  '
Set-TextValue "D13" '2.322.68'
$ws.Range("E13").Value = '  +3.61%  '
Set-TextValue "D14" '0.848'
$ws.Range("E14").Value = '  +1.07%  '
Set-TextValue "D15" '21.86'
$ws.Range("E15").Value = '  +1.65%  '
Set-TextValue "D16" '5.43'
$ws.Range("E16").Value = '  +3.09%  '
Set-TextValue "D17" '2.024.06'
$ws.Range("E17").Value = '  +3.65%  '
Set-TextValue "D18" '37.156.89'
$ws.Range("E18").Value = '  +1.90%  '
Set-TextValue "D19" '70.31'
$ws.Range("E19").Value = '  +1.53%  '
Set-TextValue "D20" '0.0₃0861'
$ws.Range("E20").Value = '  +0.86%  '
Set-TextValue "D21" '5.22'
$ws.Range("E21").Value = '  +2.73%  '
Set-TextValue "D22" '230.46'
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  +4.45%  '
$ws.Range("E25").Value = '  -0.76%  '
Set-TextValue "D26" '9.36'
$ws.Range("E26").Value = '  +2.10%  '
Set-TextValue "D27" '163.56'
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("E28").Value = '  -3.07%  '
Set-TextValue "D29" '19.75'
$ws.Range("E29").Value = '  +2.29%  '
$ws.Range("E30").Value = '  +6.01%  '
Set-TextValue "D31" '0.121'
$ws.Range("E31").Value = '  +0.74%  '
Set-TextValue "D32" '0.0671'
$ws.Range("E32").Value = '  +9.54%  '
Set-TextValue "D33" '4.76'
$ws.Range("E33").Value = '  -0.23%  '
Set-TextValue "D34" '2.50'
$ws.Range("E34").Value = '  +10.55%  '
Set-TextValue "D35" '4.45'
$ws.Range("E35").Value = '  -0.22%  '
Set-TextValue "D36" '3.61'
$ws.Range("E36").Value = '  +5.08%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("E38").Value = '  +1.73%  '
Set-TextValue "D39" '5.32'
$ws.Range("E39").Value = '  -2.31%  '
$ws.Range("E40").Value = '  +3.36%  '
Set-TextValue "D41" '0.0972'
$ws.Range("E41").Value = '  +0.35%  '
$ws.Range("E42").Value = '  +2.53%  '
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D43" '1.18'
$ws.Range("E43").Value = '  +1.38%  '
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue "D44" '16.73'
$ws.Range("E44").Value = '  +5.49%  '
Set-TextValue "D45" '91.14'
$ws.Range("E45").Value = '  +2.63%  '
Set-TextValue "D46" '1.377.65'
$ws.Range("E46").Value = '  +1.10%  '
$ws.Range("E47").Value = '  +2.52%  '
Set-TextValue "D48" '7.42'
$ws.Range("E48").Value = '  +3.68%  '
$ws.Range("E49").Value = '  +13.91%  '
Set-TextValue "D50" '2.88'
$ws.Range("E50").Value = '  +1.70%  '
Set-TextValue "D51" '46.03'
$ws.Range("E51").Value = '  +1.75%  '
